$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("02_38커뮤니케이션(최근일자기준)")

# Insert a new row at row 3 (pushes existing rows 3-21 down to 4-22)
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "신영스팩10호"
$ws.Range("B3").Value = "2024.01.22~01.23"
$ws.Range("C3").Value = "2,000~2,000"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 9150
$ws.Range("F3").Value = "신영증권"

# Remove the last row (was row 21 "NH스팩30호", now pushed to row 22)
$ws.Rows.Item(22).Delete()
